# Weekly data refresh: a new "Albahaca" price-sheet record (dated 2021-11-11,
# serial 44511) arrived and was inserted as the new row 43, pushing the
# previously-existing rows 43-53 down to rows 44-54 (each keeping its own
# data, just shifted down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 43, shifting rows 43:53 down to 44:54.
$ws.Rows(43).Insert()

# Populate the newly inserted row 43 with the latest weekly record.
$ws.Range("A43").Value = 8
$ws.Range("B43").Value = "Terminal La Palmera de La Serena"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 44511
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = 100112052
$ws.Range("G43").Value = "Albahaca"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 760
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 3500
$ws.Range("M43").Value = 3250
$ws.Range("N43").Value = "`$/paquete"
$ws.Range("O43").Value = "Región de Arica y Parinacota"
$ws.Range("P43").Value = 3250
$ws.Range("Q43").Value = 1
$ws.Range("R43").Value = "Hortaliza"
